$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the "004921978 / ELAINE / 8.08" row (row 261) and
# populate it with the "004462543 / RODOLFO / 8.83" record that the diff
# relocates earlier in the sheet.
$ws.Rows.Item(261).Insert()
$ws.Cells.Item(261, 1).Value = "'004462543"
$ws.Cells.Item(261, 2).Value = "RODOLFO"
$ws.Cells.Item(261, 3).Value = 8.83

# Remove the old trio of rows (RODOLFO -91.17 / CRISTINA -569.11 /
# DILSON -4280.85) that used to sit near the end of the data, right before
# the blank row + "Filtros aplicados" footer. After the insert above, these
# shifted down from 311-313 to 312-314.
$ws.Range("A312:A314").EntireRow.Delete()
